$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 21744082
$ws.Range("I76").Value = 25005024
$ws.Range("J76").Value = 4464.3335
$ws.Range("K76").Value = 25005024
$ws.Range("L76").Value = 4464.3335
$ws.Range("M76").Value = -25004709
$ws.Range("N76").Value = -5094.3335

$ws.Range("H79").Value = 21744082
$ws.Range("I79").Value = 25005024
$ws.Range("J79").Value = 4464.3335
$ws.Range("K79").Value = 25005024
$ws.Range("L79").Value = 4464.3335
$ws.Range("M79").Value = -25003932
$ws.Range("N79").Value = -6648.3335

$ws.Range("H100").Value = 2638.4
$ws.Range("I100").Value = 2321.6667
$ws.Range("J100").Value = 2774.1428
$ws.Range("K100").Value = 2321.6667
$ws.Range("L100").Value = 2774.1428
$ws.Range("M100").Value = -1780.6667
$ws.Range("N100").Value = -3856.1428

$ws.Range("H101").Value = 700.3
$ws.Range("J101").Value = 693
$ws.Range("L101").Value = 2079
$ws.Range("N101").Value = -5323

$ws.Range("H132").Value = 2333.4546
$ws.Range("I132").Value = 1668.7646
$ws.Range("K132").Value = 5006.293799999999
$ws.Range("M132").Value = -2476.293799999999

$ws.Range("H137").Value = 968688.75
$ws.Range("I137").Value = 1970
$ws.Range("J137").Value = 2902126.2
$ws.Range("K137").Value = 5910
$ws.Range("L137").Value = 8706378.600000001
$ws.Range("M137").Value = -3360
$ws.Range("N137").Value = -8711478.600000001

$ws.Range("H138").Value = 2207.868
$ws.Range("I138").Value = 1217
$ws.Range("J138").Value = 2761.5881
$ws.Range("K138").Value = 3651
$ws.Range("L138").Value = 8284.764299999999
$ws.Range("M138").Value = 1489
$ws.Range("N138").Value = -18564.7643

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1205.4814
$ws.Range("I2").Value = 1064.5625
$ws.Range("K2").Value = 1064.5625
$ws.Range("M2").Value = -951.5625

$ws.Range("H61").Value = 73791
$ws.Range("I61").Value = 1772.25
$ws.Range("K61").Value = 1772.25
$ws.Range("M61").Value = -1560.25

$ws.Range("H74").Value = 27024.324
$ws.Range("I74").Value = 44754.477
$ws.Range("K74").Value = 44754.477
$ws.Range("M74").Value = -43880.477

$ws.Range("H77").Value = 27024.324
$ws.Range("I77").Value = 44754.477
$ws.Range("K77").Value = 223772.385
$ws.Range("M77").Value = -219404.385

$ws.Range("H97").Value = 477.11765
$ws.Range("I97").Value = 519.2727
$ws.Range("J97").Value = 399.83334
$ws.Range("K97").Value = 519.2727
$ws.Range("L97").Value = 399.83334
$ws.Range("M97").Value = -23.27269999999999
$ws.Range("N97").Value = -1391.83334

$ws.Range("H116").Value = 1205.4814
$ws.Range("I116").Value = 1064.5625
$ws.Range("K116").Value = 1064.5625
$ws.Range("M116").Value = 1229.4375

$ws.Range("H122").Value = 3244.4856
$ws.Range("I122").Value = 3435.2415
$ws.Range("J122").Value = 2322.5
$ws.Range("K122").Value = 10305.7245
$ws.Range("L122").Value = 6967.5
$ws.Range("M122").Value = -7855.7245
$ws.Range("N122").Value = -11867.5

$ws.Range("H136").Value = 73791
$ws.Range("I136").Value = 1772.25
$ws.Range("K136").Value = 5316.75
$ws.Range("M136").Value = -2766.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1205.4814
$ws.Range("I3").Value = 1064.5625
$ws.Range("K3").Value = 1064.5625
$ws.Range("M3").Value = -950.5625

$ws.Range("H82").Value = 17953.5
$ws.Range("I82").Value = 11944.2
$ws.Range("J82").Value = 48000
$ws.Range("K82").Value = 11944.2
$ws.Range("L82").Value = 48000
$ws.Range("M82").Value = -11561.2
$ws.Range("N82").Value = -48766

$ws.Range("H85").Value = 17953.5
$ws.Range("I85").Value = 11944.2
$ws.Range("J85").Value = 48000
$ws.Range("K85").Value = 11944.2
$ws.Range("L85").Value = 48000
$ws.Range("M85").Value = -10618.2
$ws.Range("N85").Value = -50652

$ws.Range("H99").Value = 5384801.5
$ws.Range("I99").Value = 211762
$ws.Range("J99").Value = 31250000
$ws.Range("K99").Value = 211762
$ws.Range("L99").Value = 31250000
$ws.Range("M99").Value = -210264
$ws.Range("N99").Value = -31252996

$ws.Range("H134").Value = 1645.5122
$ws.Range("I134").Value = 1182.8064
$ws.Range("K134").Value = 3548.4192
$ws.Range("M134").Value = -1013.4192

$ws.Range("H140").Value = 124380.5
$ws.Range("J140").Value = 70720.57000000001
$ws.Range("L140").Value = 70720.57000000001
$ws.Range("N140").Value = -81080.57000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1936.0834
$ws.Range("J16").Value = 2499.75
$ws.Range("L16").Value = 2499.75
$ws.Range("N16").Value = -3073.75

$ws.Range("H103").Value = 27675
$ws.Range("I103").Value = 1512
$ws.Range("K103").Value = 1512
$ws.Range("M103").Value = -340

$ws.Range("H113").Value = 1936.0834
$ws.Range("J113").Value = 2499.75
$ws.Range("L113").Value = 2499.75
$ws.Range("N113").Value = -6839.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 825

$ws.Range("H60").Value = 867.59576
$ws.Range("I60").Value = 909.6
$ws.Range("K60").Value = 2728.8
$ws.Range("M60").Value = -2477.8

$ws.Range("H61").Value = 137.4
$ws.Range("I61").Value = 137.4
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 412.2
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -197.2
$ws.Range("N61").ClearContents()

$ws.Range("H119").Value = 20070.2
$ws.Range("I119").Value = 1116.6666
$ws.Range("J119").Value = 48500.5
$ws.Range("K119").Value = 3349.9998
$ws.Range("L119").Value = 145501.5
$ws.Range("M119").Value = 1488.0002
$ws.Range("N119").Value = -155177.5

$ws.Range("H121").Value = 1410.2413
$ws.Range("I121").Value = 1044.25
$ws.Range("J121").Value = 1549.6666
$ws.Range("K121").Value = 3132.75
$ws.Range("L121").Value = 4648.9998
$ws.Range("M121").Value = -1822.75
$ws.Range("N121").Value = -7268.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6410.6665
$ws.Range("I122").Value = 6102
$ws.Range("J122").Value = 6719.3335
$ws.Range("K122").Value = 18306
$ws.Range("L122").Value = 20158.0005
$ws.Range("M122").Value = -15856
$ws.Range("N122").Value = -25058.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 49030.8
$ws.Range("I7").Value = 27412.875
$ws.Range("K7").Value = 27412.875
$ws.Range("M7").Value = -27300.875

$ws.Range("H16").Value = 2244.9546
$ws.Range("I16").Value = 2476.2942
$ws.Range("J16").Value = 1458.4
$ws.Range("K16").Value = 2476.2942
$ws.Range("L16").Value = 1458.4
$ws.Range("M16").Value = -2306.2942
$ws.Range("N16").Value = -1798.4

$ws.Range("H82").Value = 1333.25
$ws.Range("I82").Value = 1412.8334
$ws.Range("J82").Value = 1094.5
$ws.Range("K82").Value = 1412.8334
$ws.Range("L82").Value = 1094.5
$ws.Range("M82").Value = -1051.8334
$ws.Range("N82").Value = -1816.5

$ws.Range("H85").Value = 1333.25
$ws.Range("I85").Value = 1412.8334
$ws.Range("J85").Value = 1094.5
$ws.Range("K85").Value = 1412.8334
$ws.Range("L85").Value = 1094.5
$ws.Range("M85").Value = -164.8334
$ws.Range("N85").Value = -3590.5

$ws.Range("H100").Value = 15793.923
$ws.Range("I100").Value = 20813.555
$ws.Range("K100").Value = 20813.555
$ws.Range("M100").Value = -20272.555

$ws.Range("H126").Value = 49030.8
$ws.Range("I126").Value = 27412.875
$ws.Range("K126").Value = 82238.625
$ws.Range("M126").Value = -79768.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 45919.6
$ws.Range("J82").Value = 47574.5
$ws.Range("L82").Value = 47574.5
$ws.Range("N82").Value = -48340.5

$ws.Range("H85").Value = 45919.6
$ws.Range("J85").Value = 47574.5
$ws.Range("L85").Value = 47574.5
$ws.Range("N85").Value = -50226.5

$ws.Range("H100").Value = 4762996
$ws.Range("I100").Value = 14286429
$ws.Range("J100").Value = 1280
$ws.Range("K100").Value = 28572858
$ws.Range("L100").Value = 2560
$ws.Range("M100").Value = -28572317
$ws.Range("N100").Value = -3642

$ws.Range("H132").Value = 1813470
$ws.Range("I132").Value = 1846.238
$ws.Range("K132").Value = 5538.714
$ws.Range("M132").Value = -3008.714
